$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 9: this pushes the old row 9 (and everything
# below it) down by one, turning the old "Child" controller block
# (rows 9-13) into rows 10-14.
$ws.Rows.Item(9).Insert() | Out-Null

# New row 8: controller button legend (Y / C / I) above the new
# "Parent" block.
$ws.Range("B8").Value = "Y"
$ws.Range("C8").Value = "C"
$ws.Range("D8").Value = "I"

# Row 10 (shifted-down former row 9) becomes the "Parent" header row
# instead of the duplicated "Child" header row.
$ws.Range("A10").Value = "Parent"
$ws.Range("A10").Font.Bold = $true

# Add a new row 15 with the same controller button legend (Y / C / I).
$ws.Range("B15").Value = "Y"
$ws.Range("C15").Value = "C"
$ws.Range("D15").Value = "I"

# Misc view/formatting tweaks captured in the diff (column A widened
# slightly to fit the new "Parent" label).
$ws.Columns.Item(1).ColumnWidth = 17.25
$ws.Range("G14").Select() | Out-Null
$ws.Application.ActiveWindow.Zoom = 145
